$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.254.93"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.77"
$ws.Range("E3").Value = "  +1.32%  "

$ws.Range("E4").Value = "  +0.70%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.78"
$ws.Range("E5").Value = "  +3.61%  "

$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.53"
$ws.Range("E8").Value = "  +6.83%  "

$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.131.38"
$ws.Range("E12").Value = "  +1.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.51"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.847.38"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.680"
$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.73"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.238.50"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.01"
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0797"
$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.26"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("E22").Value = "  +1.23%  "

$ws.Range("E23").Value = "  +0.68%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.79"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  +25.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.14"
$ws.Range("E27").Value = "  +4.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.75"
$ws.Range("E28").Value = "  +1.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.125"
$ws.Range("E29").Value = "  +0.70%  "

$ws.Range("E30").Value = "  +1.84%  "

$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.03"
$ws.Range("E32").Value = "  +2.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  +29.61%  "

$ws.Range("E34").Value = "  +2.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.07"
$ws.Range("E35").Value = "  +9.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.820"
$ws.Range("E36").Value = "  +17.83%  "

$ws.Range("E37").Value = "  +7.54%  "

$ws.Range("E38").Value = "  +3.87%  "

$ws.Range("E39").Value = "  +4.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.46"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.348.39"
$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.27"
$ws.Range("E42").Value = "  +3.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0604"
$ws.Range("E43").Value = "  +15.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  +2.86%  "

$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.44"
$ws.Range("E46").Value = "  +44.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.60"
$ws.Range("E48").Value = "  +5.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.046.63"
$ws.Range("E49").Value = "  +1.32%  "

$ws.Range("E50").Value = "  +3.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.42"
$ws.Range("E51").Value = "  +1.21%  "
